$d = $word.ActiveDocument

$ptOld = '1)Estatística Descritiva: População e amostra; apresentação gráfica dos dados; medidas de posição; medidas de dispersão.2)Amostragem: Amostragem aleatória simples com reposição; amostragem aleatória simples sem reposição.3)Conceitos de Probabilidade: Conceitos básicos de probabilidade; operações com eventos; probabilidade condicional; independência; Teorema de Bayes.4)Variáveis Aleatórias discretas: Caracterização de uma variável aleatória discreta; distribuições de probabilidade: Uniforme, Bernoulli, Binomial, Poisson, Geométrica, Binomial Negativa e Hipergeométrica.5)Variáveis Aleatórias contínuas: Caracterização de uma variável aleatória contínua; distribuições de probabilidade: Uniforme, Exponencial e Normal.6)Aproximações: Aproximação das distribuições Binomial e Poisson pela distribuição Normal.7)Teorema Central do Limite: Distribuição da média amostral; distribuição da proporção amostral; intervalos de confiança para a média amostral e para a proporção amostral; dimensionamento amostral.8)Conceitos de Testes de Hipóteses: Erro Tipo I e Erro Tipo II; p-valor; poder.9)Testes de Hipóteses para uma única amostra: Teste de hipótese para a média; teste de hipótese para a proporção e teste de hipótese para a variância.10)Testes de Hipóteses para comparação de duas amostras: Teste de hipótese para comparação de médias (amostras independentes e dependentes); teste de hipótese para comparação de duas proporções e teste de hipótese para comparação de variâncias.11) Análise de Variância: Estimação do modelo; tabela de análise de variância; intervalos de confiança para a diferença entre as médias; correção de Bonferroni; teste de homocedasticidade.12)Regressão Linear Simples e Regressão Linear Múltipla: Estimação do modelo; interpretação dos parâmetros; tabela de análise de variância; intervalos de confiança para os parâmetros; R^2; análise dos resíduos.'
$ptNew = '1)Estatística Descritiva: População e amostra; apresentação gráfica dos dados; medidas de posição; medidas de dispersão.^l2)Amostragem: Amostragem aleatória simples com reposição; amostragem aleatória simples sem reposição.^l3)Conceitos de Probabilidade: Conceitos básicos de probabilidade; operações com eventos; probabilidade condicional; independência; Teorema de Bayes.^l4)Variáveis Aleatórias discretas: Caracterização de uma variável aleatória discreta; distribuições de probabilidade: Uniforme, Bernoulli, Binomial, Poisson, Geométrica, Binomial Negativa e Hipergeométrica.^l5)Variáveis Aleatórias contínuas: Caracterização de uma variável aleatória contínua; distribuições de probabilidade: Uniforme, Exponencial e Normal.^l6)Aproximações: Aproximação das distribuições Binomial e Poisson pela distribuição Normal.^l7)Teorema Central do Limite: Distribuição da média amostral; distribuição da proporção amostral; intervalos de confiança para a média amostral e para a proporção amostral; dimensionamento amostral.^l8)Conceitos de Testes de Hipóteses: Erro Tipo I e Erro Tipo II; p-valor; poder.^l9)Testes de Hipóteses para uma única amostra: Teste de hipótese para a média; teste de hipótese para a proporção e teste de hipótese para a variância.^l10)Testes de Hipóteses para comparação de duas amostras: Teste de hipótese para comparação de médias (amostras independentes e dependentes); teste de hipótese para comparação de duas proporções e teste de hipótese para comparação de variâncias.^l11) Análise de Variância: Estimação do modelo; tabela de análise de variância; intervalos de confiança para a diferença entre as médias; correção de Bonferroni; teste de homocedasticidade.^l12)Regressão Linear Simples e Regressão Linear Múltipla: Estimação do modelo; interpretação dos parâmetros; tabela de análise de variância; intervalos de confiança para os parâmetros; R^2; análise dos resíduos.'
$d.Content.Find.Execute($ptOld, $true, $false, $false, $false, $false, $true, 1, $false, $ptNew, 2)

$enOld = '1)Descriptive Statistics: Population and sample; graphical presentation of data in statistics; measures of central tendency position and dispersion.2)Sampling methods: Simple random sampling with replacement and simple random sampling without replacement.3)Introduction to probability: Probability concepts; events probability; conditional probability; independence; Bayes Theorem.4)Discrete Random Variables: Discrete variables characterization; probability distributions: Uniform, Bernoulli, Binomial, Poisson, Geometric, Negative Binomial and Hipergeometric.5)Continuous Random Variables: Continuous variables characterization; probability distributions: Uniform, Exponential and Normal.6)Approximations: Approximation of Binomial and Poisson distributions by Normal distribution.7)Central Limit Theorem: Distribution of sample mean; distribution of sample proportion; confidence intervals for means and proportion estimated; sample sizing.8)Hypothesis test concept:  Type I Error and Type II Error; p-value; power.9)Hypothesis test for a single sample: Hypothesis test for mean; hypothesis test for proportion and hypothesis test for variance.10)Hypothesis test for two samples comparison: hypothesis test for two means comparison (dependents and independents samples); hypothesis test for two proportions comparison and hypothesis test for two variances comparison.11)Analysis of variance: Model estimation; Analysis of Variance table; confidence intervals for means difference; Bonferroni correction; homoscedasticity test.12)Simple linear regression and Multiple linear regression: Model estimation; parameters interpretation; Analysis of Variance table; confidence intervals for the parameters; R^2; residuals analysis.'
$enNew = '1)Descriptive Statistics: Population and sample; graphical presentation of data in statistics; measures of central tendency position and dispersion.^l2)Sampling methods: Simple random sampling with replacement and simple random sampling without replacement.^l3)Introduction to probability: Probability concepts; events probability; conditional probability; independence; Bayes Theorem.^l4)Discrete Random Variables: Discrete variables characterization; probability distributions: Uniform, Bernoulli, Binomial, Poisson, Geometric, Negative Binomial and Hipergeometric.^l5)Continuous Random Variables: Continuous variables characterization; probability distributions: Uniform, Exponential and Normal.^l6)Approximations: Approximation of Binomial and Poisson distributions by Normal distribution.^l7)Central Limit Theorem: Distribution of sample mean; distribution of sample proportion; confidence intervals for means and proportion estimated; sample sizing.^l8)Hypothesis test concept:  Type I Error and Type II Error; p-value; power.^l9)Hypothesis test for a single sample: Hypothesis test for mean; hypothesis test for proportion and hypothesis test for variance.^l10)Hypothesis test for two samples comparison: hypothesis test for two means comparison (dependents and independents samples); hypothesis test for two proportions comparison and hypothesis test for two variances comparison.^l11)Analysis of variance: Model estimation; Analysis of Variance table; confidence intervals for means difference; Bonferroni correction; homoscedasticity test.^l12)Simple linear regression and Multiple linear regression: Model estimation; parameters interpretation; Analysis of Variance table; confidence intervals for the parameters; R^2; residuals analysis.'
$d.Content.Find.Execute($enOld, $true, $false, $false, $false, $false, $true, 1, $false, $enNew, 2)
